# Update cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest scraped values, per the scheduled GitHub Actions refresh.
#
# D-column price strings can look numeric (e.g. "211.33"), so we force the
# cell to Text format before assigning, then clear the formatting back to
# the sheet's default (unstyled) look so only the displayed text changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.531.87"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.618.97"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.83"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0887"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.847.82"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.28%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.615.37"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.51%  "
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.551"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.52"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.526.98"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.74"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0718"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.56"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.94"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("E24").Value = "  +6.76%  "
$ws.Range("E25").Value = "  -0.37%  "
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("E27").Value = "  -1.77%  "
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.57"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.03%  "
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.443.58"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("E34").Value = "  -3.00%  "
$ws.Range("E35").Value = "  -3.13%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.945"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +4.18%  "
$ws.Range("E38").Value = "  -1.89%  "
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("E40").Value = "  -1.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.33"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.93%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  -2.00%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.40"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.28%  "
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.758.67"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.69"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.24"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("E51").Value = "  +0.37%  "
